# Update NATMI LR-pair TPM results (Wnt1-Fzd9) with newly recomputed values.
# Target cluster labels (column D) are unchanged; only the numeric
# expression/specificity metrics (columns G-T, plus K/L for the
# Resolving-Mac rows) are refreshed with the new TPM-derived figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.131499
$ws.Range("H2").Value = 0.394497
$ws.Range("I2").Value = 0.3654391092296077
$ws.Range("J2").Value = 0.3654391092296077
$ws.Range("M2").Value = 0.8794496666666666
$ws.Range("N2").Value = 2.638349
$ws.Range("O2").Value = 0.3488427963707166
$ws.Range("P2").Value = 0.3488427963707166
$ws.Range("Q2").Value = 0.115646751717
$ws.Range("R2").Value = 1.040820765453
$ws.Range("S2").Value = 0.1274808007668801
$ws.Range("T2").Value = 0.1274808007668801
$ws.Range("G3").Value = 0.131499
$ws.Range("H3").Value = 0.394497
$ws.Range("I3").Value = 0.3654391092296077
$ws.Range("J3").Value = 0.3654391092296077
$ws.Range("O3").Value = 0.2822103394539786
$ws.Range("P3").Value = 0.2822103394539786
$ws.Range("Q3").Value = 0.093557067534
$ws.Range("R3").Value = 0.842013607806
$ws.Range("S3").Value = 0.1031306950654472
$ws.Range("T3").Value = 0.1031306950654471
$ws.Range("G4").Value = 0.131499
$ws.Range("H4").Value = 0.394497
$ws.Range("I4").Value = 0.3654391092296077
$ws.Range("J4").Value = 0.3654391092296077
$ws.Range("M4").Value = 0.8404543333333333
$ws.Range("N4").Value = 2.521363
$ws.Range("O4").Value = 0.3333748945214069
$ws.Range("P4").Value = 0.3333748945214068
$ws.Range("Q4").Value = 0.110518904379
$ws.Range("R4").Value = 0.994670139411
$ws.Range("S4").Value = 0.1218282244934174
$ws.Range("T4").Value = 0.1218282244934173
$ws.Range("G5").Value = 0.131499
$ws.Range("H5").Value = 0.394497
$ws.Range("I5").Value = 0.3654391092296077
$ws.Range("J5").Value = 0.3654391092296077
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08967866666666667
$ws.Range("N5").Value = 0.269036
$ws.Range("O5").Value = 0.035571969653898
$ws.Range("P5").Value = 0.03557196965389799
$ws.Range("Q5").Value = 0.011792654988
$ws.Range("R5").Value = 0.106133894892
$ws.Range("S5").Value = 0.01299938890386312
$ws.Range("T5").Value = 0.01299938890386312
$ws.Range("G6").Value = 0.2283393333333333
$ws.Range("H6").Value = 0.685018
$ws.Range("I6").Value = 0.6345608907703922
$ws.Range("J6").Value = 0.6345608907703922
$ws.Range("M6").Value = 0.8794496666666666
$ws.Range("N6").Value = 2.638349
$ws.Range("O6").Value = 0.3488427963707166
$ws.Range("P6").Value = 0.3488427963707166
$ws.Range("Q6").Value = 0.2008129505868889
$ws.Range("R6").Value = 1.807316555282
$ws.Range("S6").Value = 0.2213619956038365
$ws.Range("T6").Value = 0.2213619956038365
$ws.Range("G7").Value = 0.2283393333333333
$ws.Range("H7").Value = 0.685018
$ws.Range("I7").Value = 0.6345608907703922
$ws.Range("J7").Value = 0.6345608907703922
$ws.Range("O7").Value = 0.2822103394539786
$ws.Range("P7").Value = 0.2822103394539786
$ws.Range("Q7").Value = 0.1624556721293333
$ws.Range("R7").Value = 1.462101049164
$ws.Range("S7").Value = 0.1790796443885314
$ws.Range("T7").Value = 0.1790796443885314
$ws.Range("G8").Value = 0.2283393333333333
$ws.Range("H8").Value = 0.685018
$ws.Range("I8").Value = 0.6345608907703922
$ws.Range("J8").Value = 0.6345608907703922
$ws.Range("M8").Value = 0.8404543333333333
$ws.Range("N8").Value = 2.521363
$ws.Range("O8").Value = 0.3333748945214069
$ws.Range("P8").Value = 0.3333748945214068
$ws.Range("Q8").Value = 0.1919087821704444
$ws.Range("R8").Value = 1.727179039534
$ws.Range("S8").Value = 0.2115466700279895
$ws.Range("T8").Value = 0.2115466700279895
$ws.Range("G9").Value = 0.2283393333333333
$ws.Range("H9").Value = 0.685018
$ws.Range("I9").Value = 0.6345608907703922
$ws.Range("J9").Value = 0.6345608907703922
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.08967866666666667
$ws.Range("N9").Value = 0.269036
$ws.Range("O9").Value = 0.035571969653898
$ws.Range("P9").Value = 0.03557196965389799
$ws.Range("Q9").Value = 0.02047716696088889
$ws.Range("R9").Value = 0.184294502648
$ws.Range("S9").Value = 0.02257258075003488
$ws.Range("T9").Value = 0.02257258075003487
